$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column D ("Customer", parallel to "Design"/"Part")
$ws.Range("D2").Value = "Customer"
$ws.Range("D2").HorizontalAlignment = -4108  # xlCenter

# Mark "Name" row (row 21) as also applicable to Customer
$ws.Range("D21").Value = "X"
$ws.Range("D21").HorizontalAlignment = -4108  # xlCenter

# New attribute rows for the Customer part, rows 25-32
$newRows = @(
    "Customer ID",
    "Customer Name",
    "Customer City",
    "Customer State",
    "Color_1",
    "Color_2",
    "Color_3",
    "School_Mascot"
)

$r = 25
foreach ($name in $newRows) {
    $ws.Cells.Item($r, 1).Value = $name

    $ws.Cells.Item($r, 4).Value = "X"
    $ws.Cells.Item($r, 4).HorizontalAlignment = -4108  # xlCenter

    $r = $r + 1
}

# Update view: scroll so row 6 is top-left, and active cell / selection is C28
$excel.ActiveWindow.ScrollRow = 6
$ws.Range("C28").Select()
